$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 78
$ws.Range("D2").Value = 42
$ws.Range("E2").Value = 32

$ws.Range("C3").Value = 99.59999999999999
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 39

$ws.Range("C4").Value = 60.3
$ws.Range("D4").Value = 48
$ws.Range("E4").Value = 28

$ws.Range("C5").Value = 81.59999999999999
$ws.Range("D5").Value = 42
$ws.Range("E5").Value = 34

$ws.Range("C6").Value = 66.7
$ws.Range("D6").Value = 45
$ws.Range("E6").Value = 30
